# Update yesterday sales table: fix a handful of Item Name / UOM pairs that
# had drifted onto the wrong row. Swap the Item Name + UOM values between
# each affected row pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-ItemUom {
    param($rowA, $rowB)

    $nameA = $ws.Range("D$rowA").Value2
    $uomA  = $ws.Range("E$rowA").Value2
    $nameB = $ws.Range("D$rowB").Value2
    $uomB  = $ws.Range("E$rowB").Value2

    $ws.Range("D$rowA").Value = $nameB
    $ws.Range("E$rowA").Value = $uomB
    $ws.Range("D$rowB").Value = $nameA
    $ws.Range("E$rowB").Value = $uomA
}

Swap-ItemUom 7 9
Swap-ItemUom 11 12
Swap-ItemUom 14 15
Swap-ItemUom 17 19
Swap-ItemUom 25 26
Swap-ItemUom 27 28
